$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DateBeg/DateEnd columns (D:E) store dates as plain text (e.g. "1995-06-01"),
# not as real Excel dates. Force those specific cells to the Text number format
# *before* writing their new values so Excel does not silently reinterpret the
# strings as date serials.
$textRanges = @("D6:D7", "D10", "D12:D30", "D32:D46", "E6:E30", "E32:E46")
foreach ($rng in $textRanges) {
    $ws.Range($rng).NumberFormat = "@"
}

# Per-cell updates (product reordering/renumbering + refreshed coverage dates,
# cross-checked against each dataset's manifest file). Only cells whose value
# actually changes are touched.
$updates = @(
    @(6, 1, 6),
    @(6, 2, "Surface Rad Budget - CCI/C3S"),
    @(6, 4, "1995-06-01"),
    @(6, 5, "2022-06-30"),
    @(7, 1, 4),
    @(7, 2, "Surface Rad Budget - CMSAF CLARA-A2"),
    @(7, 3, "Clouds"),
    @(7, 4, "1982-01-01"),
    @(7, 5, "2022-12-31"),
    @(8, 1, 5),
    @(8, 2, "Surface Rad Budget - CMSAF CLARA-A3"),
    @(8, 3, "Clouds"),
    @(8, 5, "2024-06-30"),
    @(9, 1, 10),
    @(9, 2, "Earth Rad Budget - C3S RMIB TotSolarIrrad"),
    @(9, 5, "2025-01-19"),
    @(10, 1, 9),
    @(10, 2, "Earth Rad Budget - CCI/C3S"),
    @(10, 4, "1995-06-01"),
    @(10, 5, "2022-06-30"),
    @(11, 1, 11),
    @(11, 2, "Earth Rad Budget - CMSAF CLARA-A3"),
    @(11, 5, "2024-06-30"),
    @(12, 1, 7),
    @(12, 2, "Earth Rad Budget - NASA CERES EBAF"),
    @(12, 3, "Earth Radiation Budget"),
    @(12, 4, "2000-03-01"),
    @(12, 5, "2024-07-31"),
    @(13, 1, 8),
    @(13, 2, "Earth Rad Budget - NOAA/NCEI HIRS"),
    @(13, 3, "Earth Radiation Budget"),
    @(13, 4, "1979-01-01"),
    @(13, 5, "2025-01-01"),
    @(14, 1, 15),
    @(14, 2, "Precipitation"),
    @(14, 4, "1979-01-01"),
    @(14, 5, "2024-03-31"),
    @(15, 1, 17),
    @(15, 2, "Precipitation_GIRAFE"),
    @(15, 3, "Precipitation"),
    @(15, 4, "2002-01-01"),
    @(15, 5, "2022-12-31"),
    @(16, 1, 16),
    @(16, 2, "Precipitation_microwave"),
    @(16, 3, "Precipitation"),
    @(16, 4, "2000-01-01"),
    @(16, 5, "2017-12-31"),
    @(17, 1, 14),
    @(17, 2, "Surface Rad Budget - CCI/C3S"),
    @(17, 3, "Surface Radiation Budget"),
    @(17, 4, "1995-06-01"),
    @(17, 5, "2022-06-30"),
    @(18, 1, 12),
    @(18, 2, "Surface Rad Budget - CMSAF CLARA-A2"),
    @(18, 3, "Surface Radiation Budget"),
    @(18, 4, "1982-01-01"),
    @(18, 5, "2022-12-31"),
    @(19, 1, 13),
    @(19, 2, "Surface Rad Budget - CMSAF CLARA-A3"),
    @(19, 3, "Surface Radiation Budget"),
    @(19, 4, "1979-01-01"),
    @(19, 5, "2024-06-30"),
    @(20, 2, "Total Column Water Vapour (HOAPS)"),
    @(20, 4, "1988-01-31"),
    @(20, 5, "2020-12-31"),
    @(21, 1, 19),
    @(21, 2, "Total Column Water Vapour (MERIS/SSMI)"),
    @(21, 3, "Upper-air Water Vapour"),
    @(21, 4, "2002-05-01"),
    @(21, 5, "2017-12-31"),
    @(21, 6, "Atmospheric Physics"),
    @(22, 1, 20),
    @(22, 2, "Tropospheric Humidity Profiles (RO)"),
    @(22, 3, "Upper-air Water Vapour"),
    @(22, 4, "2006-12-01"),
    @(22, 5, "2024-06-30"),
    @(22, 6, "Atmospheric Physics"),
    @(23, 1, 21),
    @(23, 2, "Upper Tropospheric Humidity"),
    @(23, 3, "Upper-air Water Vapour"),
    @(23, 4, "1999-01-01"),
    @(23, 5, "2021-02-28"),
    @(23, 6, "Atmospheric Physics"),
    @(24, 1, 26),
    @(24, 2, "Glaciers elevation and mass change data"),
    @(24, 3, "Glaciers"),
    @(24, 4, "1975-04-01"),
    @(24, 5, "2021-09-30"),
    @(25, 1, 27),
    @(25, 2, "Randolph Glacier Inventory for the year 2000"),
    @(25, 3, "Glaciers"),
    @(25, 4, "1990-01-01"),
    @(25, 5, "2010-12-31"),
    @(26, 1, 23),
    @(26, 2, "Ice Sheet Gravimetric Mass Balance"),
    @(26, 4, "2002-04-16"),
    @(26, 5, "2022-12-17"),
    @(27, 1, 24),
    @(27, 2, "Ice Sheet Surface Elevation Change (Antarctica)"),
    @(27, 3, "Ice Sheets"),
    @(27, 4, "1994-11-01"),
    @(27, 5, "2020-06-01"),
    @(27, 6, "Cryosphere"),
    @(28, 1, 25),
    @(28, 2, "Ice Sheet Surface Elevation Change (Greenland)"),
    @(28, 3, "Ice Sheets"),
    @(28, 4, "1992-01-01"),
    @(28, 5, "2024-01-01"),
    @(28, 6, "Cryosphere"),
    @(29, 1, 22),
    @(29, 2, "Ice Sheet Velocity (Greenland)"),
    @(29, 3, "Ice Sheets"),
    @(29, 4, "2018-10-01"),
    @(29, 5, "2021-09-30"),
    @(29, 6, "Cryosphere"),
    @(30, 1, 36),
    @(30, 2, "Surface Albedo 10-daily"),
    @(30, 3, "Albedo"),
    @(30, 4, "1981-09-20"),
    @(30, 5, "2020-06-30"),
    @(31, 1, 32),
    @(31, 2, "FAPAR"),
    @(31, 3, "FAPAR"),
    @(32, 1, 33),
    @(32, 2, "Fire Burned Areas"),
    @(32, 3, "Fire"),
    @(32, 4, "2001-01-01"),
    @(32, 5, "2022-12-01"),
    @(33, 1, 34),
    @(33, 2, "Fire Radiative Power"),
    @(33, 3, "Fire"),
    @(33, 4, "2020-01-01"),
    @(33, 5, "2024-02-29"),
    @(33, 6, "Land Biosphere"),
    @(34, 1, 31),
    @(34, 2, "LAI"),
    @(34, 3, "LAI"),
    @(34, 4, "1981-09-20"),
    @(34, 5, "2020-06-30"),
    @(34, 6, "Land Biosphere"),
    @(35, 1, 35),
    @(35, 2, "Land Cover"),
    @(35, 3, "Land Cover"),
    @(35, 4, "1992-01-01"),
    @(35, 5, "2022-12-31"),
    @(35, 6, "Land Biosphere"),
    @(36, 1, 28),
    @(36, 2, "Lake Surface Temperature"),
    @(36, 3, "Lakes"),
    @(36, 4, "1995-06-01"),
    @(36, 5, "2023-12-31"),
    @(36, 6, "Land Hydrology"),
    @(37, 1, 29),
    @(37, 2, "Lake Water Level"),
    @(37, 3, "Lakes"),
    @(37, 4, "1992-09-26"),
    @(37, 5, "2023-12-30"),
    @(37, 6, "Land Hydrology"),
    @(38, 1, 30),
    @(38, 2, "Soil Moisture"),
    @(38, 3, "Soil Moisture"),
    @(38, 4, "1978-11-01"),
    @(38, 5, "2024-12-31"),
    @(38, 6, "Land Hydrology"),
    @(39, 1, 37),
    @(39, 2, "Ocean Colour"),
    @(39, 3, "Ocean Colour"),
    @(39, 4, "1997-09-04"),
    @(39, 5, "2024-09-30"),
    @(40, 1, 41),
    @(40, 2, "SST"),
    @(40, 3, "SST"),
    @(40, 4, "1981-08-24"),
    @(40, 5, "2022-12-31"),
    @(41, 1, 42),
    @(41, 2, "SST (ESA CCI GMPE)"),
    @(41, 3, "SST"),
    @(41, 4, "1981-09-01"),
    @(41, 5, "2016-12-31"),
    @(42, 1, 38),
    @(42, 2, "Sea Ice Concentration"),
    @(42, 3, "Sea Ice"),
    @(42, 4, "1978-10-25"),
    @(42, 5, "2025-01-01"),
    @(43, 1, 39),
    @(43, 2, "Sea Ice Edge and Type"),
    @(43, 3, "Sea Ice"),
    @(43, 4, "1978-10-25"),
    @(43, 5, "2025-01-01"),
    @(44, 1, 40),
    @(44, 2, "Sea Ice Thickness"),
    @(44, 3, "Sea Ice"),
    @(44, 4, "2002-10-01"),
    @(44, 5, "2024-04-30"),
    @(44, 6, "Ocean"),
    @(45, 1, 43),
    @(45, 2, "Sea Level"),
    @(45, 3, "Sea Level"),
    @(45, 4, "1993-01-01"),
    @(45, 5, "2023-06-07"),
    @(45, 6, "Ocean"),
    @(46, 1, 44),
    @(46, 2, "Surface Geostrophic Currents"),
    @(46, 3, "Surface Currents"),
    @(46, 4, "1993-01-01"),
    @(46, 5, "2023-06-07"),
    @(46, 6, "Ocean")
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# The sheet originally only went down to row 43; rows 44-46 are brand new.
# Give column A in those new rows the same border/alignment styling already
# used by every other row's "#" column (copy the formatting from A2).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A44:A46").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

